# Auto-generated Excel COM-interop script to apply profit recalculation updates
# across the Halicarnassus_Profits workbook sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 96.85714
$ws.Range("I9").Value = 111.666664
$ws.Range("K9").Value = 111.666664
$ws.Range("M9").Value = 57.333336
# Row 11
$ws.Range("H11").Value = 13.941176
$ws.Range("I11").Value = 13.941176
$ws.Range("K11").Value = 13.941176
$ws.Range("M11").Value = 126.058824
# Row 19
$ws.Range("H19").Value = 494.2
$ws.Range("I19").Value = 377.375
$ws.Range("J19").Value = 627.7143
$ws.Range("K19").Value = 377.375
$ws.Range("L19").Value = 627.7143
$ws.Range("M19").Value = -202.375
$ws.Range("N19").Value = -977.7143
# Row 29
$ws.Range("H29").Value = 3693.625
$ws.Range("I29").Value = 1849.6666
$ws.Range("J29").Value = 4800
$ws.Range("K29").Value = 5548.9998
$ws.Range("L29").Value = 14400
$ws.Range("M29").Value = -5267.9998
$ws.Range("N29").Value = -14962
# Row 38
$ws.Range("H38").Value = 1220.4286
$ws.Range("J38").Value = 3500
$ws.Range("L38").Value = 10500
$ws.Range("N38").Value = -11244
# Row 58
$ws.Range("H58").Value = 2229.75
$ws.Range("J58").Value = 3312.75
$ws.Range("L58").Value = 9938.25
$ws.Range("N58").Value = -10238.25
# Row 138
$ws.Range("H138").Value = 7899
$ws.Range("J138").Value = 11250
$ws.Range("L138").Value = 33750
$ws.Range("N138").Value = -44030

$ws = $wb.Worksheets.Item("ARM")
# Row 6
$ws.Range("H6").Value = 12860428
$ws.Range("I6").Value = 10913272
$ws.Range("K6").Value = 10913272
$ws.Range("M6").Value = -10913099
# Row 19
$ws.Range("H19").Value = 9049.333000000001
$ws.Range("I19").Value = 9049.333000000001
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 9049.333000000001
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -8820.333000000001
$ws.Range("N19").ClearContents()
# Row 44
$ws.Range("H44").Value = 11620.5
$ws.Range("J44").Value = 11620.5
$ws.Range("L44").Value = 11620.5
$ws.Range("N44").Value = -12596.5
# Row 55
$ws.Range("H55").Value = 36197.8
$ws.Range("J55").Value = 41497.25
$ws.Range("L55").Value = 41497.25
$ws.Range("N55").Value = -42127.25

$ws = $wb.Worksheets.Item("BSM")
# Row 35
$ws.Range("H35").Value = 77777
$ws.Range("J35").Value = 77777
$ws.Range("L35").Value = 77777
$ws.Range("N35").Value = -78397
# Row 134
$ws.Range("H134").Value = 9670
$ws.Range("I134").Value = 3428.75
$ws.Range("K134").Value = 10286.25
$ws.Range("M134").Value = -7751.25

$ws = $wb.Worksheets.Item("CRP")
# Row 68
$ws.Range("H68").Value = 84682.5
$ws.Range("J68").Value = 84682.5
$ws.Range("L68").Value = 84682.5
$ws.Range("N68").Value = -86180.5
# Row 71
$ws.Range("H71").Value = 84682.5
$ws.Range("J71").Value = 84682.5
$ws.Range("L71").Value = 254047.5
$ws.Range("N71").Value = -261535.5
# Row 99
$ws.Range("H99").Value = 4328.3335
$ws.Range("I99").Value = 4328.3335
$ws.Range("K99").Value = 4328.3335
$ws.Range("M99").Value = -2830.3335
# Row 126
$ws.Range("H126").Value = 4328.3335
$ws.Range("I126").Value = 4328.3335
$ws.Range("K126").Value = 12985.0005
$ws.Range("M126").Value = -10515.0005

$ws = $wb.Worksheets.Item("CUL")
# Row 34
$ws.Range("H34").Value = 1966.3572
$ws.Range("J34").Value = 2252.9167
$ws.Range("L34").Value = 6758.750100000001
$ws.Range("N34").Value = -6926.750100000001
# Row 39
$ws.Range("H39").Value = 5332.222
$ws.Range("J39").Value = 5332.222
$ws.Range("L39").Value = 15996.666
$ws.Range("N39").Value = -16584.666
# Row 55
$ws.Range("H55").Value = 2493
$ws.Range("I55").Value = 425.125
$ws.Range("J55").Value = 3674.6428
$ws.Range("K55").Value = 1275.375
$ws.Range("L55").Value = 11023.9284
$ws.Range("M55").Value = -1098.375
$ws.Range("N55").Value = -11377.9284

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 251.73914
$ws.Range("I2").Value = 179.5
$ws.Range("K2").Value = 179.5
$ws.Range("M2").Value = -66.5
# Row 10
$ws.Range("H10").Value = 1032.5
$ws.Range("I10").Value = 1500
$ws.Range("J10").Value = 876.6667
$ws.Range("K10").Value = 1500
$ws.Range("L10").Value = 876.6667
$ws.Range("M10").Value = -1331
$ws.Range("N10").Value = -1214.6667
# Row 19
$ws.Range("H19").Value = 2466.6667
$ws.Range("J19").Value = 2500
$ws.Range("L19").Value = 2500
$ws.Range("N19").Value = -3076
# Row 55
$ws.Range("H55").Value = 5754.6665
$ws.Range("I55").Value = 6507.5
$ws.Range("J55").Value = 4249
$ws.Range("K55").Value = 6507.5
$ws.Range("L55").Value = 4249
$ws.Range("M55").Value = -6180.5
$ws.Range("N55").Value = -4903

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3260
$ws.Range("I7").Value = 2100
$ws.Range("J7").Value = 5000
$ws.Range("K7").Value = 2100
$ws.Range("L7").Value = 5000
$ws.Range("M7").Value = -1988
$ws.Range("N7").Value = -5224
# Row 11
$ws.Range("H11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()
# Row 40
$ws.Range("H40").Value = 7840.5
$ws.Range("I40").Value = 7840.5
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 7840.5
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -7704.5
$ws.Range("N40").ClearContents()
# Row 46
$ws.Range("H46").Value = 4827.9375
$ws.Range("I46").Value = 666.6667
$ws.Range("J46").Value = 5788.231
$ws.Range("K46").Value = 666.6667
$ws.Range("L46").Value = 5788.231
$ws.Range("M46").Value = -478.6667
$ws.Range("N46").Value = -6164.231
# Row 55
$ws.Range("H55").Value = 1393.6
$ws.Range("I55").Value = 945.8182
$ws.Range("K55").Value = 945.8182
$ws.Range("M55").Value = -772.8182
# Row 58
$ws.Range("H58").Value = 1500
$ws.Range("I58").Value = 1500
$ws.Range("K58").Value = 1500
$ws.Range("M58").Value = -1240
# Row 93
$ws.Range("H93").Value = 1875.5454
$ws.Range("I93").Value = 1842.125
$ws.Range("J93").Value = 1964.6666
$ws.Range("K93").Value = 1842.125
$ws.Range("L93").Value = 1964.6666
$ws.Range("M93").Value = -594.125
$ws.Range("N93").Value = -4460.6666
# Row 126
$ws.Range("H126").Value = 3260
$ws.Range("I126").Value = 2100
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 6300
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -3830
$ws.Range("N126").Value = -19940

$ws = $wb.Worksheets.Item("WVR")
# Row 6
$ws.Range("H6").Value = 549.75
$ws.Range("I6").Value = 533
$ws.Range("K6").Value = 533
$ws.Range("M6").Value = -418
# Row 7
$ws.Range("H7").Value = 1518.75
$ws.Range("J7").Value = 2000
$ws.Range("L7").Value = 2000
$ws.Range("N7").Value = -2226
# Row 81
$ws.Range("H81").Value = 1100
$ws.Range("I81").Value = 1100
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 2200
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -1139
$ws.Range("N81").ClearContents()
# Row 84
$ws.Range("H84").Value = 1100
$ws.Range("I84").Value = 1100
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 11000
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -5696
$ws.Range("N84").ClearContents()
# Row 96
$ws.Range("H96").Value = 1756.3334
$ws.Range("I96").Value = 1715.2858
$ws.Range("K96").Value = 1715.2858
$ws.Range("M96").Value = -342.2858000000001
